# Fruta / hortaliza, semanal
# Insert a new weekly record at row 211 (pushing the existing rows 211..263
# down to 212..264) and populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 211; this shifts rows 211-263 down
# to 212-264 and extends the used range to R264 automatically.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new weekly record. The
# non-varying descriptive columns (A, B, C, E, F, G, H, I, R) repeat the same
# values used throughout the whole table.
$ws.Cells.Item(211, 1).Value = 7
$ws.Cells.Item(211, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(211, 3).Value = "Ñuble"
$ws.Cells.Item(211, 4).Value = 44855
$ws.Cells.Item(211, 5).Value = 16
$ws.Cells.Item(211, 6).Value = 100112043
$ws.Cells.Item(211, 7).Value = "Pepino ensalada"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 60
$ws.Cells.Item(211, 11).Value = 20000
$ws.Cells.Item(211, 12).Value = 20000
$ws.Cells.Item(211, 13).Value = 20000
$ws.Cells.Item(211, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(211, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(211, 16).Value = 333
$ws.Cells.Item(211, 17).Value = 60
$ws.Cells.Item(211, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# the "Fecha" column.
$ws.Cells.Item(211, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
